$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.186522666666667
$ws.Range("H2").Value = 3.559568
$ws.Range("I2").Value = 0.06400371352898657
$ws.Range("J2").Value = 0.06400371352898658
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 2.405218290947555
$ws.Range("R2").Value = 21.646964618528
$ws.Range("S2").Value = 0.0004221867076807287
$ws.Range("T2").Value = 0.0004221867076807288

# Row 3
$ws.Range("G3").Value = 1.186522666666667
$ws.Range("H3").Value = 3.559568
$ws.Range("I3").Value = 0.06400371352898657
$ws.Range("J3").Value = 0.06400371352898658
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 304.2769810155875
$ws.Range("R3").Value = 2738.492829140288
$ws.Range("S3").Value = 0.05340957921428161
$ws.Range("T3").Value = 0.05340957921428161

# Row 4
$ws.Range("G4").Value = 1.186522666666667
$ws.Range("H4").Value = 3.559568
$ws.Range("I4").Value = 0.06400371352898657
$ws.Range("J4").Value = 0.06400371352898658
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 57.95008225952179
$ws.Range("R4").Value = 521.5507403356961
$ws.Range("S4").Value = 0.01017194760702424
$ws.Range("T4").Value = 0.01017194760702424

# Row 5
$ws.Range("I5").Value = 0.599039184070822
$ws.Range("J5").Value = 0.599039184070822
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 22.51150633422089
$ws.Range("R5").Value = 202.603557007988
$ws.Range("S5").Value = 0.003951432923967324
$ws.Range("T5").Value = 0.003951432923967324

# Row 6
$ws.Range("I6").Value = 0.599039184070822
$ws.Range("J6").Value = 0.599039184070822
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.4998839753196394
$ws.Range("T6").Value = 0.4998839753196394

# Row 7
$ws.Range("I7").Value = 0.599039184070822
$ws.Range("J7").Value = 0.599039184070822
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 542.3805601195186
$ws.Range("R7").Value = 4881.425041075667
$ws.Range("S7").Value = 0.09520377582721537
$ws.Range("T7").Value = 0.09520377582721534

# Row 8
$ws.Range("G8").Value = 6.246625666666667
$ws.Range("H8").Value = 18.739877
$ws.Range("I8").Value = 0.3369571024001913
$ws.Range("J8").Value = 0.3369571024001914
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 12.66263067049355
$ws.Range("R8").Value = 113.963676034442
$ws.Range("S8").Value = 0.002222664933770562
$ws.Range("T8").Value = 0.002222664933770562

# Row 9
$ws.Range("G9").Value = 6.246625666666667
$ws.Range("H9").Value = 18.739877
$ws.Range("I9").Value = 0.3369571024001913
$ws.Range("J9").Value = 0.3369571024001914
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 1601.911579765703
$ws.Range("R9").Value = 14417.20421789133
$ws.Range("S9").Value = 0.2811827011304164
$ws.Range("T9").Value = 0.2811827011304164

# Row 10
$ws.Range("G10").Value = 6.246625666666667
$ws.Range("H10").Value = 18.739877
$ws.Range("I10").Value = 0.3369571024001913
$ws.Range("J10").Value = 0.3369571024001914
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 305.0868570802188
$ws.Range("R10").Value = 2745.781713721969
$ws.Range("S10").Value = 0.05355173633600441
$ws.Range("T10").Value = 0.05355173633600439

